$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to make room for the new first forecast vector,
# shifting the existing rows 2:52 down to 3:53.
$ws.Rows("2:2").Insert()

# The insert copies row-1 (header) formatting onto the new row by default.
# Re-establish the original look: column A keeps the date-number style (same
# as used throughout the rest of column A), columns B:E stay unstyled.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("B2:E2").ClearFormats()

# Full refreshed data set (date serial, y_0, y_0_forecast, y_1, y_1_forecast)
# for every row, rows 2 through 53.
$data = @(
    @(39400, 2007, 1.75539628881467, 2008, 0.2337905658324813),
    @(39583, 2008, 2.071001150990881, 2009, 5.125398993052044),
    @(39765, 2008, 2.213911448916162, 2009, 3.386383090739953),
    @(39948, 2009, 0.8787496612562951, 2010, 1.013823151053095),
    @(40130, 2009, 2.533533936850563, 2010, 0.984293482975751),
    @(40310, 2010, 2.723861837742825, 2011, 4.356912452939476),
    @(40494, 2010, 2.088987486264915, 2011, 3.612753212925401),
    @(40676, 2011, 2.431458940167008, 2012, 5.259925231829898),
    @(40862, 2011, 1.212544822741002, 2012, 2.158838189283174),
    @(41044, 2012, 1.447930496829541, 2013, 0.7814198158142105),
    @(41228, 2012, 1.196776590518644, 2013, 1.194058515117336),
    @(41409, 2013, 0.3494148569448852, 2014, -0.5632089058212553),
    @(41592, 2013, 0.4712609263772594, 2014, 1.409662779709819),
    @(41774, 2014, 0.6742451383205061, 2015, 1.713290556413605),
    @(41957, 2014, 0.8783377572271434, 2015, 2.372074663906587),
    @(42137, 2015, 2.178094576990031, 2016, 2.743122680804988),
    @(42321, 2015, 2.29066283401107, 2016, 4.595879021798321),
    @(42503, 2016, 2.623600596229347, 2017, 1.891220645890002),
    @(42689, 2016, 4.109890522944348, 2017, 4.034919509273061),
    @(42867, 2017, 1.715791310593229, 2018, 1.687339605296523),
    @(43053, 2017, 1.336316831462692, 2018, 0.02883756256675252),
    @(43145, 2018, 1.862609889357314, 2019, 2.016108144125295),
    @(43235, 2018, 0.3317798769387315, 2019, -2.079848588862154),
    @(43326, 2018, 1.369357173039498, 2019, 2.147322685428343),
    @(43418, 2018, 1.197912858979611, 2019, 0.9262553939922924),
    @(43510, 2019, 4.206171608278875, 2020, 6.673057635354218),
    @(43600, 2019, 0.8934982674867697, 2020, -1.194610791899997),
    @(43691, 2019, 1.749048192229496, 2020, 2.097071260293459),
    @(43783, 2019, 1.727537197898665, 2020, 2.928189816005666),
    @(43875, 2020, 2.306826470345347, 2021, 1.405199465409468),
    @(43966, 2020, 1.979074033580819, 2021, 0.8024032015999882),
    @(44068, 2020, 2.980209378995857, 2021, 2.026374749120596),
    @(44159, 2020, 3.647228437274408, 2021, 3.673004547855219),
    @(44251, 2021, 2.128447415063373, 2022, 2.829537440100038),
    @(44341, 2021, 2.552476296061434, 2022, 3.086122033237126),
    @(44432, 2021, 1.954146674711188, 2022, -1.871731962523027),
    @(44525, 2021, 2.777797690741424, 2022, 1.579011422502852),
    @(44617, 2022, 1.014265466411501, 2023, 0.00000000000006661338147750939),
    @(44706, 2022, 0.3901728183783204, 2023, 1.906002353653125),
    @(44798, 2022, 2.69102598245059, 2023, 3.932998599265303),
    @(44890, 2022, 0.6994919452575576, 2023, -2.087978868409623),
    @(44981, 2023, -0.6123001687638907, 2024, -1.985049937500016),
    @(45071, 2023, -1.906744368254853, 2024, 9.556389850000446),
    @(45163, 2023, -1.669605379075589, 2024, -0.4516362914963601),
    @(45254, 2023, -1.432689847121871, 2024, 0.1172571542027212),
    @(45345, 2024, 0.5893837960974757, 2025, -0.3994003998999518),
    @(45436, 2024, 1.152760694685062, 2025, 0.6444718444275521),
    @(45534, 2024, 1.780300968358017, 2025, 0.8104204674762139),
    @(45618, 2024, 2.033479419175133, 2025, 1.317145539573517),
    @(45713, 2025, 2.506151357112452, 2026, 0),
    @(45800, 2025, 2.287212358310953, 2026, 1.801540135156521),
    @(45891, 2025, 2.481068287768839, 2026, 1.642273543406181)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
